# Insert "the " before "highest number" in the first highlight bullet.
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Bangladesh reported highest number of annual dengue cases and deaths in 2023 ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Bangladesh reported the highest number of annual dengue cases and deaths in 2023 ",
    2)
